$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing content and rebuild the report table so the shared-string
# table gets regenerated in top-to-bottom first-use order (matches a fresh
# report run).
$ws.Cells.Clear()

$ws.Range("A1").Value = "Number"
$ws.Range("B1").Value = "Studentenname"
$ws.Range("C1").Value = "MissingCourses"

$data = @(
    @(1, "Student84", 2),
    @(2, "Student45", 1),
    @(3, "Student79", 1),
    @(4, "Student5", 1),
    @(5, "Student113", 1),
    @(6, "Student86", 2),
    @(7, "Student60", 1),
    @(8, "Student2", 1),
    @(9, "Student134", 1),
    @(10, "Student10", 2),
    @(11, "Student110", 2),
    @(12, "Student114", 2),
    @(13, "Student90", 2),
    @(14, "Student22", 1),
    @(15, "Student130", 2),
    @(16, "Student118", 1),
    @(17, "Student137", 1),
    @(18, "Student111", 1),
    @(19, "Student150", 1),
    @(20, "Student87", 1),
    @(21, "Student148", 1),
    @(22, "Student109", 2),
    @(23, "Student49", 1),
    @(24, "Student139", 1),
    @(25, "Student36", 1),
    @(26, "Student50", 1),
    @(27, "Student141", 1),
    @(28, "Student125", 2),
    @(29, "Student37", 1),
    @(30, "Student77", 2),
    @(31, "Student33", 1),
    @(32, "Student115", 1),
    @(33, "Student132", 1),
    @(34, "Student51", 1),
    @(35, "Student92", 1),
    @(36, "Student129", 1),
    @(37, "Student53", 1),
    @(38, "Student35", 1),
    @(39, "Student131", 1),
    @(40, "Student18", 1),
    @(41, "Student81", 1),
    @(42, "Student102", 2),
    @(43, "Student94", 1),
    @(44, "Student103", 1),
    @(45, "Student143", 1),
    @(46, "Student70", 2),
    @(47, "Student121", 1),
    @(48, "Student127", 1),
    @(49, "Student48", 1),
    @(50, "Student107", 1),
    @(51, "Student74", 1),
    @(52, "Student126", 1),
    @(53, "Student93", 1),
    @(54, "Student42", 1),
    @(55, "Student61", 1),
    @(56, "Student38", 1),
    @(57, "Student23", 1),
    @(58, "Student106", 1),
    @(59, "Student78", 1),
    @(60, "Student122", 1),
    @(61, "Student76", 1),
    @(62, "Student34", 2),
    @(63, "Student116", 1),
    @(64, "Student138", 1),
    @(65, "Student98", 1),
    @(66, "Student62", 1),
    @(67, "Student80", 1),
    @(68, "Student100", 1),
    @(69, "Student142", 1),
    @(70, "Student54", 1)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
